# Update cryptocurrency price/volume data per latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.228.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.587.05'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.37%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '316.31'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.68%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.578'
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.541'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '35.85'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0814'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.53'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.983.07'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.29%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.108'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.549.15'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.847'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.256.18'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.77%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.86'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +3.11%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.81%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0964'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '69.60'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '254.77'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '27.24'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '10.34'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '5.87'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '156.83'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.37%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.42'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0809'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('E37').Value = '  -1.39%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.48'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +8.12%  '
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '22.71'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.97'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.53%  '
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.26'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.015.72'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.98'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '83.27'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.76%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '75.88'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.47%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.814.93'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('E51').Value = '  +2.55%  '
